# Append "for COVID-19 in Malaysia" to the research-question text on the
# title slide (Slide 1) and on the research-question recap slide (Slide 5).
#
# Both edits replace only the trailing run of text (".../in each cluster?")
# with "... in each cluster for COVID-19 in Malaysia?" so that existing run
# formatting / other runs in the paragraph are left untouched.

$p = $ppt.ActivePresentation

function Update-ClusterQuestion($Shape, $OldSubstring, $NewSubstring) {
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldSubstring)
    if ($idx -ge 0) {
        $startPos = $idx + 1
        $len = $OldSubstring.Length
        $sub = $tr.Characters($startPos, $len)
        $sub.Text = $NewSubstring
    }
}

# Slide 1 - "Title 1": full research question sentence ends in "...cluster?"
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(1)
Update-ClusterQuestion $title1 "Is there a correlation between the number of tests conducted & total number of cases reported in each cluster?" "Is there a correlation between the number of tests conducted & total number of cases reported in each cluster for COVID-19 in Malaysia?"

# Slide 5 - "Title 4": same question, last run is just "in each cluster?"
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(4)
Update-ClusterQuestion $title5 "in each cluster?" "in each cluster for COVID-19 in Malaysia?"
